$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# --- Update "Description" (G) steps text for row 10 ---
# Inserts an extra "validate4;" + "wait(2);" checkpoint after the VT200-0438 run, and
# renumbers the final validation (for VT200-0439) from validate4 to validate5.
$g10 = "wait(3);`nvalidate1;`nlink_Click(intent_test_link);`nvalidate2;`nSelectTestToRun(VT200_0438_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(2);`nvalidate4;`nwait(2);`nSelectTestToRun(VT200_0439_string);`nClickRunTest(runtest_top_xpath);`nvalidate5;`nClickRunTest(runtest_bottom_xpath);`nwait(2);`npress_Key(Home);`nlaunch_App_Device(com.rhomobile.compliancetest_js/com.rhomobile.rhodes.RhodesActivity);`nCheckUITextContains(stoplistening);"
$ws.Range("G10").Value = $g10

# --- Update "Expected Behaviour" (H) text for row 10 (VT200-0438 + VT200-0439 combined) ---
# Now includes validate4 (the broadcast alert) plus a new validate5 block asserting VT200-0439.
$h10 = "validate1`n{`nvalidate_PageTitle=Compliance JS specs`n};`nvalidate2`n{`nvalidate_PageTitle=Intent JS Test`n};`nvalidate3`n{`nvalidate_OldText_Exists=VT200-0438`n};`nvalidate4`n{`nvalidate_Alert={`"data`":{`"myData`":`"This is broad cast data 3!`"},`"appName`":`"com.rhomobile.compliancetest_js`",`"action`":`"com.rhomobile.BROADCAST`",`"intentType`":`"broadcast`"};`nvalidate5`n{`nvalidate_OldText_Exists=VT200-0439`n};`n"
$ws.Range("H10").Value = $h10

# --- Update "Expected Behaviour" (H) text for row 9 (VT200-0438 broadcast alert) ---
# Adds a trailing semicolon after the validate_Alert JSON payload's closing brace.
$h9 = "validate1`n{`nvalidate_PageTitle=Compliance JS specs`n};`nvalidate2`n{`nvalidate_PageTitle=Intent JS Test`n};`nvalidate3`n{`nvalidate_OldText_Exists=VT200-0438`n};`nvalidate4`n{`nvalidate_Alert={`"data`":{`"myData`":`"This is broad cast data 3!`"},`"appName`":`"com.rhomobile.compliancetest_js`",`"action`":`"com.rhomobile.BROADCAST`",`"intentType`":`"broadcast`"};`n"
$ws.Range("H9").Value = $h9

# Row 10 grew taller to fit the extra validation lines (screenshots re-captured at Nexus7 size).
$ws.Rows.Item(10).RowHeight = 306.75

# Selection moved from D1 to E2.
$ws.Range("E2").Select()
